$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above row 3 for the "Sniper/Long Range/Overwatch" section header ---
$ws.Rows("3:3").Insert()

# New row 3 content: header spanning A3:B3
$ws.Range("A3").Value = "Sniper/Long Range/Overwatch"
$ws.Range("A3:B3").Merge()

# B3 (the merged-away cell) drops back to a plain, non-bold/underlined, non-top-aligned look
$ws.Range("B3").Font.Underline = $false
$ws.Range("B3").Font.Bold = $false
$ws.Range("B3").VerticalAlignment = -4107

# --- 2. Update ability names that were renamed (WOTC_APA_Sprint -> WOTC_APA_Phase) ---
$ws.Range("F6").Value = "WOTC_APA_Phase"
$ws.Range("F8").Value = "WOTC_APA_Phase"

# --- 3. Fill in newly added columns for the Sentry/Default row (now row 10) ---
$ws.Range("E10").Value = "WOTC_APA_EverVigilant"
$ws.Range("F10").Value = "WOTC_APA_WeaponsHot"
$ws.Range("G10").Value = "WOTC_APA_CombatAwareness"
$ws.Range("H10").Value = "WOTC_APA_Sentinel"

# --- 4. Re-center the column header row (row 2, columns C:J) ---
$ws.Range("C2:J2").HorizontalAlignment = -4108

# --- 5. Column width adjustments ---
$ws.Columns("B").ColumnWidth = 18.7109375
$ws.Columns("G").ColumnWidth = 30.28515625

# --- 6. View adjustments ---
$ws.Range("E20").Select()
